$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set cell B3 to "y" (matches A3 = "y_col")
$ws.Range("B3").Value = "y"
